$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D ("D:E"), shifting the existing
# quarterly columns (old D:K) right to F:M. This makes room for the two
# newest reporting quarters.
$ws.Columns("D:E").Insert()

# Copy the number formatting from column F (which now holds what used to be
# column D) onto the freshly inserted D:E columns so the new cells pick up
# the same date / number styles as the rest of the table instead of the
# generic default.
$ws.Range("F7:F102").Copy()
$ws.Range("D7:E102").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("D7").Value = 43461
$ws.Range("E7").Value = 43370
$ws.Range("D8").Value = 137400
$ws.Range("E8").Value = 110100
$ws.Range("D9").Value = 34100
$ws.Range("E9").Value = 33200
$ws.Range("D10").Value = 103300
$ws.Range("E10").Value = 76900
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 1500
$ws.Range("E14").Value = 3200
$ws.Range("H14").Value = -191600
$ws.Range("I14").Value = 0
$ws.Range("D15").Value = 10400
$ws.Range("E15").Value = 10000
$ws.Range("D17").Value = 78100
$ws.Range("E17").Value = 71000
$ws.Range("H17").Value = -121100
$ws.Range("D18").Value = 59300
$ws.Range("E18").Value = 39100
$ws.Range("H18").Value = 261800
$ws.Range("D20").Value = 1500
$ws.Range("E20").Value = 700
$ws.Range("H20").Value = 9400
$ws.Range("I20").Value = 3400
$ws.Range("D21").Value = 71200
$ws.Range("E21").Value = 49800
$ws.Range("H21").Value = 280600
$ws.Range("I21").Value = 63100
$ws.Range("D22").Value = 13100
$ws.Range("E22").Value = 14400
$ws.Range("D23").Value = 47700
$ws.Range("E23").Value = 25400
$ws.Range("H23").Value = 248800
$ws.Range("I23").Value = 37400
$ws.Range("D24").Value = 6800
$ws.Range("E24").Value = -300
$ws.Range("H24").Value = 179300
$ws.Range("I24").Value = 1200
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 40900
$ws.Range("E26").Value = 25700
$ws.Range("H26").Value = 69500
$ws.Range("I26").Value = 36200
$ws.Range("D27").Value = 16300
$ws.Range("E27").Value = 11200
$ws.Range("H27").Value = 40700
$ws.Range("I27").Value = 13700
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("E29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = -1500
$ws.Range("E32").Value = -700
$ws.Range("H32").Value = -9400
$ws.Range("I32").Value = -3400
$ws.Range("D33").Value = 16300
$ws.Range("E33").Value = 11200
$ws.Range("H33").Value = -77500
$ws.Range("I33").Value = 13700
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = 16300
$ws.Range("E35").Value = 11200
$ws.Range("H35").Value = -77500
$ws.Range("I35").Value = 13700
$ws.Range("D38").Value = 43461
$ws.Range("E38").Value = 43370
$ws.Range("D41").Value = 41400
$ws.Range("E41").Value = 27300
$ws.Range("D42").Value = 24000
$ws.Range("E42").Value = 29200
$ws.Range("D43").Value = 161600
$ws.Range("E43").Value = 122400
$ws.Range("D44").Value = 0
$ws.Range("E44").Value = 0
$ws.Range("D45").Value = 3900
$ws.Range("E45").Value = 4400
$ws.Range("D46").Value = 230900
$ws.Range("E46").Value = 183300
$ws.Range("D47").Value = 13200
$ws.Range("E47").Value = 17500
$ws.Range("D48").Value = 33600
$ws.Range("E48").Value = 32800
$ws.Range("D49").Value = 684500
$ws.Range("E49").Value = 699500
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 179600
$ws.Range("E52").Value = 186900
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 1141800
$ws.Range("E54").Value = 1120000
$ws.Range("D57").Value = 18000
$ws.Range("E57").Value = 17600
$ws.Range("D58").Value = 2700
$ws.Range("E58").Value = 2700
$ws.Range("D59").Value = 89800
$ws.Range("E59").Value = 73600
$ws.Range("D60").Value = 110500
$ws.Range("E60").Value = 93900
$ws.Range("D61").Value = 920900
$ws.Range("E61").Value = 914800
$ws.Range("D62").Value = 199600
$ws.Range("E62").Value = 201700
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 1509800
$ws.Range("E66").Value = 1492200
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = -153600
$ws.Range("E72").Value = -156800
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = -368000
$ws.Range("E76").Value = -372200
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43461
$ws.Range("E80").Value = 43370
$ws.Range("D81").Value = 16300
$ws.Range("E81").Value = 11200
$ws.Range("H81").Value = -77500
$ws.Range("I81").Value = 13700
$ws.Range("D83").Value = 10400
$ws.Range("E83").Value = 10000
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 31500
$ws.Range("E89").Value = 52600
$ws.Range("D91").Value = -4200
$ws.Range("E91").Value = -2800
$ws.Range("H91").Value = -3900
$ws.Range("I91").Value = -1900
$ws.Range("J91").Value = -2900
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = 3900
$ws.Range("E94").Value = -10900
$ws.Range("D96").Value = -13200
$ws.Range("E96").Value = -13100
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = -21300
$ws.Range("E100").Value = -48800
$ws.Range("D101").Value = 0
$ws.Range("E101").Value = 0
$ws.Range("D102").Value = 14100
$ws.Range("E102").Value = -7100

$excel.CutCopyMode = $false
